{"js": "// Commit: 02/02/2018 HARISH CHICK IN\n//\n// 1) Normalize the \"TUE JAN 30 ... 2018\" timestamp paragraph so its text\n//    lives in a single run (the diff merges two adjacent <w:r> runs into\n//    one).\n// 2) Append a brand-new purchase-detail block (THU FEB 01 timestamp +\n//    Person Name / Bill number / Item Name / ... / Amount balance rows)\n//    right after the final existing transaction in the log, followed by\n//    the two trailing spacer paragraphs the diff adds.\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Part 1: merge the \"TUE JAN 30\" / \" 11:20:52 PST 2018\" runs into one.\n// ---------------------------------------------------------------------\nconst tsResults = body.search(\"TUE JAN 30*2018\", { matchWildcards: true });\ntsResults.load(\"text\");\nawait context.sync();\n\nif (tsResults.items.length === 0) {\n  throw new Error('Could not find the \"TUE JAN 30 ... 2018\" paragraph.');\n}\n\nconst tsRange = tsResults.items[0];\nconst tsText = tsRange.text; // \"TUE JAN 30 11:20:52 PST 2018\"\n// Re-writing the whole matched range collapses it back down to a single\n// run, which is exactly what the diff does (same font, so formatting is\n// unaffected).\ntsRange.insertText(tsText, Word.InsertLocation.replace);\nawait context.sync();\n\n// ---------------------------------------------------------------------\n// Part 2: append the new 02/02/2018 purchase block.\n// ---------------------------------------------------------------------\nconst COURIER = \"Courier New\";\n\n// Build the compact OOXML for one \"Label <tabs> - Value\" row, matching\n// the run layout Word produces: the label in its own run, each filler\n// tab in its own run, and the final tab sharing a run with the value.\nfunction rPr(bold) {\n  return (\n    '<w:rPr><w:rFonts w:ascii=\"' + COURIER + '\" w:hAnsi=\"' + COURIER +\n    '\" w:cs=\"' + COURIER + '\"/>' + (bold ? \"<w:b/>\" : \"\") + \"</w:rPr>\"\n  );\n}\n\nfunction runText(text, bold, preserveSpace) {\n  const sp = preserveSpace ? ' xml:space=\"preserve\"' : \"\";\n  return \"<w:r>\" + rPr(bold) + \"<w:t\" + sp + \">\" + text + \"</w:t></w:r>\";\n}\n\nfunction runTab(bold) {\n  return \"<w:r>\" + rPr(bold) + \"<w:tab/></w:r>\";\n}\n\nfunction runTabText(text, bold) {\n  return \"<w:r>\" + rPr(bold) + \"<w:tab/><w:t>\" + text + \"</w:t></w:r>\";\n}\n\nfunction para(innerRuns, bold) {\n  return (\n    \"<w:p><w:pPr><w:pStyle w:val=\\\"PlainText\\\"/>\" + rPr(bold) + \"</w:pPr>\" +\n    innerRuns + \"</w:p>\"\n  );\n}\n\nfunction emptyPara(bold) {\n  return \"<w:p><w:pPr><w:pStyle w:val=\\\"PlainText\\\"/>\" + rPr(bold) + \"</w:pPr></w:p>\";\n}\n\nfunction row(label, value, tabCount, bold) {\n  let runs = runText(label, bold);\n  for (let i = 0; i < tabCount - 1; i++) {\n    runs += runTab(bold);\n  }\n  runs += runTabText(value, bold);\n  return para(runs, bold);\n}\n\nlet newBlock = \"\";\nnewBlock += emptyPara(true);\nnewBlock += para(runText(\"THU FEB 01\", false) + runText(\" 12:38:02 PST 2018\", false, true), false);\nnewBlock += row(\"Person Name\", \"- TSV\", 4, false);\nnewBlock += row(\"Bill number\", \"- 2897\", 4, false);\nnewBlock += para(runText(\"---------------------------------------------------------------\", false), false);\nnewBlock += row(\"Item Name\", \"- CHOWCHOW\", 4, false);\nnewBlock += row(\"Number of Pockets\", \"- 1\", 3, false);\nnewBlock += row(\"Number of KGs\", \"- 88\", 3, false);\nnewBlock += row(\"Rate\", \"- 5\", 5, false);\nnewBlock += row(\"Total Price\", \"- 440.0\", 4, false);\nnewBlock += row(\"Amount balance\", \"- 800.0\", 3, true);\nnewBlock += emptyPara(false);\nnewBlock += emptyPara(true);\n\nconst ooxmlPackage =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  newBlock +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\n// Find the very last \"Amount balance\" row in the log (the one ending in\n// \"- 360.0\") so the new block lands at the end of the document, after\n// the final existing transaction.\nbody.paragraphs.load(\"items,text\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs.items;\nlet lastAmountBalanceIndex = -1;\nfor (let i = 0; i < paragraphs.length; i++) {\n  if (paragraphs[i].text.indexOf(\"Amount balance\") !== -1) {\n    lastAmountBalanceIndex = i;\n  }\n}\n\nif (lastAmountBalanceIndex === -1) {\n  throw new Error('Could not find an \"Amount balance\" paragraph.');\n}\n\nconst lastAmountBalancePara = paragraphs[lastAmountBalanceIndex];\n\n// Insert a throwaway placeholder paragraph right after it, then replace\n// that placeholder's whole range with our precise OOXML block. Doing the\n// insert in two steps (instead of inserting OOXML directly on a\n// collapsed \"after\" range) keeps the original \"Amount balance\" paragraph\n// untouched.\nconst placeholder = lastAmountBalancePara.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\nconst placeholderRange = placeholder.getRange(Word.RangeLocation.whole);\nplaceholderRange.insertOoxml(ooxmlPackage, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Commit: 02/02/2018 HARISH CHICK IN\n#\n# 1) Normalize the \"TUE JAN 30 ... 2018\" timestamp paragraph so its text\n#    lives in a single run (the diff merges two adjacent runs into one).\n# 2) Append a brand-new purchase-detail block (THU FEB 01 timestamp +\n#    Person Name / Bill number / Item Name / ... / Amount balance rows)\n#    right after the final existing transaction in the log, followed by\n#    the two trailing spacer paragraphs the diff adds.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Part 1: merge the \"TUE JAN 30\" / \" 11:20:52 PST 2018\" runs into one.\n# ---------------------------------------------------------------------\n$find = $d.Content.Find\n$find.Text = \"TUE JAN 30*2018\"\n$find.MatchWildcards = $true\n$found = $find.Execute()\nif (-not $found) {\n  throw 'Could not find the \"TUE JAN 30 ... 2018\" paragraph.'\n}\n\n$tsRange = $find.Parent\n$tsText = $tsRange.Text   # \"TUE JAN 30 11:20:52 PST 2018\"\n\n# Delete the two-run span and retype it as plain text; re-applying the\n# Courier New font afterwards collapses everything back down to a single\n# run sharing one <w:rPr>, matching what the diff does.\n$tsRange.Delete()\n$tsRange.InsertAfter($tsText)\n$tsRange.Font.NameAscii = \"Courier New\"\n$tsRange.Font.NameOther = \"Courier New\"\n$tsRange.Font.NameBi = \"Courier New\"\n\n# ---------------------------------------------------------------------\n# Part 2: append the new 02/02/2018 purchase block.\n# ---------------------------------------------------------------------\n\nfunction Get-RPr([bool]$bold) {\n  $b = \"\"\n  if ($bold) { $b = \"<w:b/>\" }\n  return \"<w:rPr><w:rFonts w:ascii=`\"Courier New`\" w:hAnsi=`\"Courier New`\" w:cs=`\"Courier New`\"/>$b</w:rPr>\"\n}\n\nfunction Get-RunText([string]$text, [bool]$bold, [bool]$preserve) {\n  $sp = \"\"\n  if ($preserve) { $sp = \" xml:space=`\"preserve`\"\" }\n  $rpr = Get-RPr $bold\n  return \"<w:r>$rpr<w:t$sp>$text</w:t></w:r>\"\n}\n\nfunction Get-RunTab([bool]$bold) {\n  $rpr = Get-RPr $bold\n  return \"<w:r>$rpr<w:tab/></w:r>\"\n}\n\nfunction Get-RunTabText([string]$text, [bool]$bold) {\n  $rpr = Get-RPr $bold\n  return \"<w:r>$rpr<w:tab/><w:t>$text</w:t></w:r>\"\n}\n\nfunction Get-Para([string]$runs, [bool]$bold) {\n  $rpr = Get-RPr $bold\n  return \"<w:p><w:pPr><w:pStyle w:val=`\"PlainText`\"/>$rpr</w:pPr>$runs</w:p>\"\n}\n\nfunction Get-EmptyPara([bool]$bold) {\n  $rpr = Get-RPr $bold\n  return \"<w:p><w:pPr><w:pStyle w:val=`\"PlainText`\"/>$rpr</w:pPr></w:p>\"\n}\n\nfunction Get-Row([string]$label, [string]$value, [int]$tabCount, [bool]$bold) {\n  $runs = Get-RunText $label $bold $false\n  for ($i = 0; $i -lt ($tabCount - 1); $i++) {\n    $runs += Get-RunTab $bold\n  }\n  $runs += Get-RunTabText $value $bold\n  return Get-Para $runs $bold\n}\n\n$block = \"\"\n$block += Get-EmptyPara $true\n$block += Get-Para ((Get-RunText \"THU FEB 01\" $false $false) + (Get-RunText \" 12:38:02 PST 2018\" $false $true)) $false\n$block += Get-Row \"Person Name\" \"- TSV\" 4 $false\n$block += Get-Row \"Bill number\" \"- 2897\" 4 $false\n$block += Get-Para (Get-RunText \"---------------------------------------------------------------\" $false $false) $false\n$block += Get-Row \"Item Name\" \"- CHOWCHOW\" 4 $false\n$block += Get-Row \"Number of Pockets\" \"- 1\" 3 $false\n$block += Get-Row \"Number of KGs\" \"- 88\" 3 $false\n$block += Get-Row \"Rate\" \"- 5\" 5 $false\n$block += Get-Row \"Total Price\" \"- 440.0\" 4 $false\n$block += Get-Row \"Amount balance\" \"- 800.0\" 3 $true\n$block += Get-EmptyPara $false\n$block += Get-EmptyPara $true\n\n$ooxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + $block + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>'\n\n# Find the very last \"Amount balance\" row in the log (the one ending in\n# \"- 360.0\") so the new block lands at the end of the document, after the\n# final existing transaction.\n$count = $d.Paragraphs.Count\n$lastAmountBalanceIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n  if ($d.Paragraphs.Item($i).Range.Text -like \"*Amount balance*\") {\n    $lastAmountBalanceIdx = $i\n  }\n}\nif ($lastAmountBalanceIdx -eq -1) {\n  throw 'Could not find an \"Amount balance\" paragraph.'\n}\n\n$lastAmountBalanceRange = $d.Paragraphs.Item($lastAmountBalanceIdx).Range\n\n# Insert a throwaway placeholder paragraph right after it, then replace\n# that placeholder's content with our precise OOXML block via InsertXML.\n# Doing the insert in two steps (instead of calling InsertXML directly on\n# a collapsed \"after\" range) keeps the original \"Amount balance\"\n# paragraph untouched.\n$lastAmountBalanceRange.InsertParagraphAfter()\n\n$placeholder = $d.Paragraphs.Item($lastAmountBalanceIdx + 1).Range\n[void]$placeholder.InsertXML($ooxml)\n"}
